# Updates cryptos list figures (prices + 1h volume change %) and
# swaps the dogwifhat/Bittensor row ordering, per the "Updated cryptos
# list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.453.63'
$ws.Range("E2").Value = '  -4.17%  '
$ws.Range("D3").Value = '3.336.04'
$ws.Range("E3").Value = '  -0.57%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '''573.23'
$ws.Range("E5").Value = '  -3.10%  '
$ws.Range("D6").Value = '''180.29'
$ws.Range("E6").Value = '  -5.85%  '
$ws.Range("D7").Value = '''0.622'
$ws.Range("E7").Value = '  +2.69%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  -3.69%  '
$ws.Range("D10").Value = '''6.67'
$ws.Range("D11").Value = '''0.402'
$ws.Range("E11").Value = '  -3.61%  '
$ws.Range("D12").Value = '3.918.99'
$ws.Range("E12").Value = '  -0.61%  '
$ws.Range("E13").Value = '  -1.08%  '
$ws.Range("D14").Value = '66.572.78'
$ws.Range("E14").Value = '  -4.07%  '
$ws.Range("D15").Value = '''26.72'
$ws.Range("E15").Value = '  -6.18%  '
$ws.Range("D16").Value = '''0.0000166'
$ws.Range("E16").Value = '  -2.64%  '
$ws.Range("D17").Value = '3.330.84'
$ws.Range("E17").Value = '  -0.77%  '
$ws.Range("D18").Value = '''432.89'
$ws.Range("E18").Value = '  -3.86%  '
$ws.Range("D19").Value = '''13.55'
$ws.Range("E19").Value = '  -1.52%  '
$ws.Range("D20").Value = '''5.66'
$ws.Range("E20").Value = '  -2.72%  '
$ws.Range("E21").Value = '  -3.19%  '
$ws.Range("D22").Value = '''73.40'
$ws.Range("E22").Value = '  -3.48%  '
$ws.Range("D23").Value = '''0.999'
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("D24").Value = '''0.517'
$ws.Range("E24").Value = '  -1.12%  '
$ws.Range("E25").Value = '  -4.44%  '
$ws.Range("E26").Value = '  +1.27%  '
$ws.Range("D27").Value = '''9.01'
$ws.Range("E27").Value = '  -4.65%  '
$ws.Range("E28").Value = '  +0.08%  '
$ws.Range("D29").Value = '''1.95'
$ws.Range("E29").Value = '  -3.14%  '
$ws.Range("D30").Value = '''22.78'
$ws.Range("E30").Value = '  -2.30%  '
$ws.Range("D32").Value = '''5.23'
$ws.Range("E32").Value = '  -5.80%  '
$ws.Range("D33").Value = '''6.76'
$ws.Range("E33").Value = '  -3.23%  '
$ws.Range("D34").Value = '''1.21'
$ws.Range("E34").Value = '  -5.77%  '
$ws.Range("D35").Value = '''160.11'
$ws.Range("E36").Value = '  -6.19%  '
$ws.Range("D37").Value = '''27.67'
$ws.Range("E37").Value = '  +1.66%  '
$ws.Range("E38").Value = '  -7.73%  '
$ws.Range("D39").Value = '2.809.37'
$ws.Range("E39").Value = '  +3.82%  '
$ws.Range("D40").Value = '''0.797'
$ws.Range("E40").Value = '  -1.25%  '
$ws.Range("E41").Value = '  -3.92%  '
$ws.Range("D42").Value = '''6.19'
$ws.Range("E42").Value = '  -4.39%  '
$ws.Range("D43").Value = '''40.24'
$ws.Range("E43").Value = '  -2.19%  '
$ws.Range("D44").Value = '''0.0666'
$ws.Range("E44").Value = '  -3.65%  '
$ws.Range("D45").Value = '''24.17'
$ws.Range("E45").Value = '  -4.87%  '
$ws.Range("B46").Value = 'Bittensor'
$ws.Range("C46").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D46").Value = '''324.90'
$ws.Range("E46").Value = '  -2.90%  '
$ws.Range("B47").Value = 'dogwifhat'
$ws.Range("C47").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D47").Value = '''2.32'
$ws.Range("E47").Value = '  -7.18%  '
$ws.Range("E48").Value = '  -4.51%  '
$ws.Range("E49").Value = '  +0.67%  '
$ws.Range("D51").Value = '''6.14'
$ws.Range("E51").Value = '  -2.93%  '
